$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: Incline bench press -> Incline chest press, mark Machine column, clear note
$ws.Range("B2").Value = "x"
$ws.Range("C2").Value = "Incline chest press"
$ws.Range("E2").Value = ""

# Row 7: Romanian deadlift -> freebie, clear weight and note
$ws.Range("C7").Value = "freebie"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""

# Row 8: add note
$ws.Range("E8").Value = "fix 3/30"

# Reset cursor/selection back to A1
$ws.Range("A1").Select()
